# The edit swaps the content of rows 16 and 19 for columns
# A, B, E, F, G, H, Q, R, Z, AB (all other columns are identical
# between the two rows already, so nothing else needs to change).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R", "Z", "AB")

foreach ($col in $cols) {
    $cell16 = $ws.Range($col + "16")
    $cell19 = $ws.Range($col + "19")
    $v16 = $cell16.Value2
    $v19 = $cell19.Value2
    $cell16.Value = $v19
    $cell19.Value = $v16
}
